$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.462.44'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '3.502.74'
$ws.Range('E3').Value = '  +4.20%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.40'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '666.71'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.48'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +6.24%  '
$ws.Range('E8').Value = '  +1.74%  '
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.999'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('D11').Value = '3.501.43'
$ws.Range('E11').Value = '  +4.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.55'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +11.60%  '
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('E14').Value = '  +2.38%  '
$ws.Range('D15').Value = '98.212.20'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000263'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('D17').Value = '4.155.66'
$ws.Range('E17').Value = '  +4.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.98'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.77%  '
$ws.Range('D19').Value = '3.500.86'
$ws.Range('E19').Value = '  +4.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.97'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +11.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.84'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +9.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.526'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '527.26'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.45'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000204'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.83'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +9.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '98.75'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.77'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.73%  '
$ws.Range('D29').Value = '3.690.86'
$ws.Range('E29').Value = '  +4.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.48'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +12.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.88'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +13.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.147'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.190'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.610'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +10.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.56'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +10.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.54'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.99'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('E40').Value = '  +4.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '528.36'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.938'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +10.76%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.77'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.44'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0438'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.81'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.97%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.72'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('B49').Value = 'MantraDAO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.64'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('E50').Value = '  +10.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.33'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.53%  '
